$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '61.558.97'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').Value = '3.447.02'
$ws.Range('E3').Value = '  +1.93%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'579.82"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.15%  '
$ws.Range('D6').Value = "'149.70"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +9.08%  '
$ws.Range('D7').Value = '3.448.37'
$ws.Range('E7').Value = '  +2.03%  '
$ws.Range('E9').Value = '  +0.88%  '
$ws.Range('D10').Value = "'7.80"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.39%  '
$ws.Range('E11').Value = '  +1.24%  '
$ws.Range('E12').Value = '  +1.30%  '
$ws.Range('D13').Value = '4.035.04'
$ws.Range('E13').Value = '  +1.95%  '
$ws.Range('D14').Value = "'27.92"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +6.68%  '
$ws.Range('E15').Value = '  -0.48%  '
$ws.Range('E16').Value = '  +0.80%  '
$ws.Range('D17').Value = '3.442.86'
$ws.Range('E17').Value = '  +1.68%  '
$ws.Range('D18').Value = '61.666.28'
$ws.Range('E18').Value = '  +0.92%  '
$ws.Range('E19').Value = '  +8.70%  '
$ws.Range('D20').Value = "'14.30"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.94%  '
$ws.Range('D21').Value = "'9.50"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.59%  '
$ws.Range('D22').Value = "'388.72"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +3.17%  '
$ws.Range('E23').Value = '  +2.66%  '
$ws.Range('D24').Value = '3.586.36'
$ws.Range('E24').Value = '  +1.72%  '
$ws.Range('D25').Value = "'72.99"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.46%  '
$ws.Range('E26').Value = '  +0.73%  '
$ws.Range('E27').Value = '  +0.16%  '
$ws.Range('E28').Value = '  -1.63%  '
$ws.Range('D29').Value = "'0.182"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +5.62%  '
$ws.Range('D30').Value = "'7.78"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +3.56%  '
$ws.Range('D31').Value = "'1.00"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('D32').Value = "'1.52"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -13.46%  '
$ws.Range('E33').Value = '  +1.35%  '
$ws.Range('E34').Value = '  +1.11%  '
$ws.Range('E36').Value = '  +1.40%  '
$ws.Range('D37').Value = "'5.28"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.12%  '
$ws.Range('D38').Value = "'7.06"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.76%  '
$ws.Range('E39').Value = '  +1.16%  '
$ws.Range('D40').Value = "'166.68"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.62%  '
$ws.Range('E41').Value = '  +3.97%  '
$ws.Range('E42').Value = '  +9.82%  '
$ws.Range('D43').Value = "'0.793"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +2.10%  '
$ws.Range('E44').Value = '  +2.01%  '
$ws.Range('D45').Value = "'0.999"
$ws.Range('D45').ClearFormats()
$ws.Range('E47').Value = '  -0.14%  '
$ws.Range('D48').Value = '2.609.84'
$ws.Range('E48').Value = '  +5.36%  '
$ws.Range('E49').Value = '  -2.75%  '
$ws.Range('D50').Value = "'7.04"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.55%  '
$ws.Range('D51').Value = "'23.18"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.35%  '
